# "sistemazione formattazione ppt e pdf"
# On slide 11, the small summary table ("Tabella 9") has its first data
# row (row 2: "6" / "24" / "446.91") reformatted: each cell's paragraph
# gets an explicit left alignment and the cell itself is vertically
# centered (anchor="ctr").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(8)
$tbl = $shape.Table

for ($col = 1; $col -le 3; $col++) {
    $cell = $tbl.Cell(2, $col)
    $tf = $cell.Shape.TextFrame

    # <a:pPr algn="l"/> on the cell's paragraph
    $tf.TextRange.ParagraphFormat.Alignment = 1

    # <a:tcPr anchor="ctr"/> on the cell
    $tf.VerticalAnchor = 3
}
